# Insert a new price-record row at row 576 (pushes existing rows 576..659 down
# to 577..660) and populate it with the new "Asterix / 2a (guarda)" entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(576).Insert()

$ws.Range("A576").Value2 = 7
$ws.Range("B576").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C576").Value2 = "Ñuble"
$ws.Range("D576").Value2 = 45124
$ws.Range("E576").Value2 = 16
$ws.Range("F576").Value2 = 100114001
$ws.Range("G576").Value2 = "Papa"
$ws.Range("H576").Value2 = "Asterix"
$ws.Range("I576").Value2 = "2a (guarda)"
$ws.Range("J576").Value2 = 120
$ws.Range("K576").Value2 = 16000
$ws.Range("L576").Value2 = 16000
$ws.Range("M576").Value2 = 16000
$ws.Range("N576").Value2 = "`$/saco 25 kilos"
$ws.Range("O576").Value2 = "Región de Los Lagos"
$ws.Range("P576").Value2 = 640
$ws.Range("Q576").Value2 = 25
$ws.Range("R576").Value2 = "Hortaliza"
